$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E3").Value = 13
$ws.Range("E6").Value = 3
$ws.Range("E8").Value = 9
$ws.Range("E9").Value = 11
$ws.Range("E15").Value = 114
$ws.Range("F15").Value = 52
$ws.Range("H15").Value = 52
$ws.Range("E17").Value = 71
$ws.Range("F17").Value = 26
$ws.Range("H17").Value = 26
$ws.Range("E18").Value = 63
$ws.Range("E19").Value = 30
$ws.Range("E24").Value = 15
$ws.Range("F25").Value = 5
$ws.Range("H25").Value = 5
$ws.Range("E27").Value = 7
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 3
$ws.Range("E31").Value = 2
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = 1
$ws.Range("E33").Value = 20
$ws.Range("F36").Value = 19
$ws.Range("H36").Value = 19
$ws.Range("E37").Value = 28
$ws.Range("F37").Value = 13
$ws.Range("H37").Value = 13
$ws.Range("E41").Value = 21
$ws.Range("F41").Value = 9
$ws.Range("H41").Value = 9
$ws.Range("E45").Value = 16
$ws.Range("E48").Value = 15
$ws.Range("E49").Value = 43
$ws.Range("F49").Value = 22
$ws.Range("H49").Value = 22
$ws.Range("E57").Value = 9
$ws.Range("E59").Value = 6
$ws.Range("F59").Value = 1
$ws.Range("H59").Value = 1
$ws.Range("E62").Value = 22
$ws.Range("E66").Value = 24
$ws.Range("F66").Value = 10
$ws.Range("H66").Value = 10
$ws.Range("E67").Value = 24
$ws.Range("E70").Value = 20
$ws.Range("F70").Value = 7
$ws.Range("H70").Value = 7
$ws.Range("E71").Value = 18
$ws.Range("E73").Value = 17
$ws.Range("F75").Value = 4
$ws.Range("H75").Value = 4
$ws.Range("E82").Value = 3
$ws.Range("E84").Value = 3
$ws.Range("E88").Value = 10
$ws.Range("F88").Value = 6
$ws.Range("H88").Value = 6
$ws.Range("F89").Value = 8
$ws.Range("H89").Value = 8
